# Add a new data row (row 77) to the COVID cases sheet, mirroring the
# structure of the existing rows (A..R), for 13/06/2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

$ws.Cells.Item($row, 1).Value = 76
$ws.Cells.Item($row, 2).Value = "13/06/2020"
$ws.Cells.Item($row, 3).Value = 268
$ws.Cells.Item($row, 4).Value = 9
$ws.Cells.Item($row, 5).Value = 41
$ws.Cells.Item($row, 6).Value = "176,4458022"
$ws.Cells.Item($row, 7).Value = "0,03358208955"
$ws.Cells.Item($row, 8).Value = 217
$ws.Cells.Item($row, 9).Value = 469
$ws.Cells.Item($row, 10).Value = 737
$ws.Cells.Item($row, 11).Value = 15
$ws.Cells.Item($row, 12).Value = 47
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = 43
$ws.Cells.Item($row, 15).Value = 35
$ws.Cells.Item($row, 16).Value = 16
$ws.Cells.Item($row, 17).Value = 10
$ws.Cells.Item($row, 18).Value = 11
